{"js": "// Remove the obsolete \"R22-3\" alternate-path table row (\"In the step 2,\n// system warns that volunteer with the same username is already logged\n// in.\") per the commit \"Fixed use-cases. Removed unnecessary use-cases.\"\n\nconst table = context.document.body.tables.getFirst();\nconst rows = table.rows;\nrows.load(\"items/values\");\nawait context.sync();\n\nlet targetRow = null;\nfor (const row of rows.items) {\n  const values = row.values || [];\n  const rowText = values.map((r) => r.join(\" \")).join(\" \");\n  if (rowText.indexOf(\"R22-3\") !== -1) {\n    targetRow = row;\n    break;\n  }\n}\n\nif (!targetRow) {\n  throw new Error(\"Could not locate the R22-3 row to delete.\");\n}\n\ntargetRow.delete();\nawait context.sync();\n", "ps1": "# Remove the obsolete \"R22-3\" alternate-path table row (\"In the step 2,\n# system warns that volunteer with the same username is already logged\n# in.\") per the commit \"Fixed use-cases. Removed unnecessary use-cases.\"\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$targetIndex = -1\nfor ($i = 1; $i -le $t.Rows.Count; $i++) {\n    $row = $t.Rows.Item($i)\n    if ($row.Range.Text.Contains(\"R22-3\")) {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not locate the R22-3 row to delete.\"\n}\n\n$t.Rows.Item($targetIndex).Delete()\n"}
